$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("11.01.2021", 38140, 168895, 1902, 24326, 11912, 0),
    @("10.01.2021", 38099, 167957, 1889, 24270, 11940, 0),
    @("09.01.2021", 38002, 166952, 1875, 24117, 12010, 0),
    @("08.01.2021", 37853, 165764, 1868, 23943, 12042, 0),
    @("05.01.2021", 37519, 164164, 1828, 23624, 12067, 0),
    @("03.01.2021", 37335, 163324, 1801, 23437, 12097, 0)
)

$startRow = 234
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($row, $col)
        if ($c -eq 0) {
            # Force the date-looking string to be stored as literal text,
            # matching the rest of column A, without leaving a stray style.
            $cell.NumberFormat = "@"
            $cell.Value = $rowData[$c]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $rowData[$c]
        }
    }
}
